$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-coerced to a number by Excel;
# force them to Text format first, then restore the default "Normal" style so no
# spurious formatting diff is introduced.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D15", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D28", "D29", "D30", "D32", "D33", "D38", "D39", "D40", "D41", "D43", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($c in $textCells) { $ws.Range($c).NumberFormat = "@" }

# Apply the updated values
$ws.Range("D2").Value = "26.660.42"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.643.31"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "215.44"
$ws.Range("E5").Value = "  +1.12%  "
$ws.Range("D6").Value = "0.504"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "0.252"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "19.11"
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "1.871.69"
$ws.Range("E12").Value = "  +1.06%  "
$ws.Range("D13").Value = "1.640.96"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("D15").Value = "0.531"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "64.95"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "26.662.49"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").Value = "0.0₃0744"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "216.37"
$ws.Range("E19").Value = "  +0.52%  "
$ws.Range("D20").Value = "1.01"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "4.35"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "9.47"
$ws.Range("E23").Value = "  +1.83%  "
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  +12.67%  "
$ws.Range("D25").Value = "145.44"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  +4.10%  "
$ws.Range("D29").Value = "15.70"
$ws.Range("E29").Value = "  +0.97%  "
$ws.Range("D30").Value = "0.0513"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("E31").Value = "  +1.77%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +1.40%  "
$ws.Range("D33").Value = "3.00"
$ws.Range("E33").Value = "  +2.02%  "
$ws.Range("D34").Value = "1.275.68"
$ws.Range("E34").Value = "  +4.63%  "
$ws.Range("E35").Value = "  +3.29%  "
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("D38").Value = "0.534"
$ws.Range("E38").Value = "  +6.13%  "
$ws.Range("D39").Value = "0.821"
$ws.Range("E39").Value = "  +3.05%  "
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "0.811"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "5.43"
$ws.Range("E43").Value = "  +1.19%  "
$ws.Range("D44").Value = "1.780.80"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("D45").Value = "91.74"
$ws.Range("E45").Value = "  -1.17%  "
$ws.Range("D46").Value = "59.26"
$ws.Range("E46").Value = "  +8.32%  "
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  +1.21%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.76"
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0964"
$ws.Range("E50").Value = "  +1.27%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "0.407"
$ws.Range("E51").Value = "  -0.11%  "

# Restore default style on the cells we forced to Text format
foreach ($c in $textCells) { $ws.Range($c).Style = "Normal" }
